# Update the 'paises.xlsx' country-data table (Pais sheet):
# - refresh case numbers for the countries whose stats changed in this run
# - the country list stays sorted by 'Casos totales' (column B) descending,
#   so rows are rewritten in the resulting sorted order
# - bump the 'Datos actualizados...' timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 216,8
$data[0,0]='Estados Unidos'; $data[0,1]=2160323; $data[0,2]=18099; $data[0,3]=859527; $data[0,4]=1182962; $data[0,5]=0; $data[0,6]=307; $data[0,7]=117834
$data[1,0]='Brasil'; $data[1,1]=852785; $data[1,2]=1989; $data[1,3]=437512; $data[1,4]=372436; $data[1,5]=0; $data[1,6]=46; $data[1,7]=42837
$data[2,0]='Rusia'; $data[2,1]=528964; $data[2,2]=8835; $data[2,3]=280050; $data[2,4]=241966; $data[2,5]=0; $data[2,6]=119; $data[2,7]=6948
$data[3,0]='India'; $data[3,1]=333008; $data[3,2]=11382; $data[3,3]=169689; $data[3,4]=153799; $data[3,5]=0; $data[3,6]=321; $data[3,7]=9520
$data[4,0]='Reino Unido'; $data[4,1]=295889; $data[4,2]=1514; $data[4,3]=0; $data[4,4]=0; $data[4,5]=0; $data[4,6]=36; $data[4,7]=41698
$data[5,0]='España'; $data[5,1]=291008; $data[5,2]=323; $data[5,3]=0; $data[5,4]=0; $data[5,5]=0; $data[5,6]=0; $data[5,7]=27136
$data[6,0]='Italia'; $data[6,1]=236989; $data[6,2]=338; $data[6,3]=176370; $data[6,4]=26274; $data[6,5]=0; $data[6,6]=44; $data[6,7]=34345
$data[7,0]='Peru'; $data[7,1]=225132; $data[7,2]=0; $data[7,3]=111724; $data[7,4]=106910; $data[7,5]=0; $data[7,6]=0; $data[7,7]=6498
$data[8,0]='Alemania'; $data[8,1]=187631; $data[8,2]=208; $data[8,3]=172200; $data[8,4]=6562; $data[8,5]=0; $data[8,6]=2; $data[8,7]=8869
$data[9,0]='Iran'; $data[9,1]=187427; $data[9,2]=2472; $data[9,3]=148674; $data[9,4]=29916; $data[9,5]=0; $data[9,6]=107; $data[9,7]=8837
$data[10,0]='Turquia'; $data[10,1]=178239; $data[10,2]=1562; $data[10,3]=151417; $data[10,4]=22015; $data[10,5]=0; $data[10,6]=15; $data[10,7]=4807
$data[11,0]='Chile'; $data[11,1]=174293; $data[11,2]=6938; $data[11,3]=143704; $data[11,4]=27266; $data[11,5]=0; $data[11,6]=222; $data[11,7]=3323
$data[12,0]='Francia'; $data[12,1]=157220; $data[12,2]=407; $data[12,3]=72859; $data[12,4]=54954; $data[12,5]=0; $data[12,6]=9; $data[12,7]=29407
$data[13,0]='Mexico'; $data[13,1]=142690; $data[13,2]=3494; $data[13,3]=104975; $data[13,4]=20843; $data[13,5]=0; $data[13,6]=424; $data[13,7]=16872
$data[14,0]='Pakistan'; $data[14,1]=139230; $data[14,2]=6825; $data[14,3]=51735; $data[14,4]=84863; $data[14,5]=0; $data[14,6]=81; $data[14,7]=2632
$data[15,0]='Arabia Saudita'; $data[15,1]=127541; $data[15,2]=4233; $data[15,3]=84720; $data[15,4]=41849; $data[15,5]=0; $data[15,6]=40; $data[15,7]=972
$data[16,0]='Canada'; $data[16,1]=98735; $data[16,2]=325; $data[16,3]=60241; $data[16,4]=30348; $data[16,5]=0; $data[16,6]=39; $data[16,7]=8146
$data[17,0]='Banglades'; $data[17,1]=87520; $data[17,2]=3141; $data[17,3]=18730; $data[17,4]=67619; $data[17,5]=0; $data[17,6]=32; $data[17,7]=1171
$data[18,0]='China'; $data[18,1]=83132; $data[18,2]=57; $data[18,3]=78369; $data[18,4]=129; $data[18,5]=0; $data[18,6]=0; $data[18,7]=4634
$data[19,0]='Catar'; $data[19,1]=79602; $data[19,2]=1186; $data[19,3]=56898; $data[19,4]=22631; $data[19,5]=0; $data[19,6]=3; $data[19,7]=73
$data[20,0]='Sudafrica'; $data[20,1]=70038; $data[20,2]=4302; $data[20,3]=38531; $data[20,4]=30027; $data[20,5]=0; $data[20,6]=57; $data[20,7]=1480
$data[21,0]='Belgica'; $data[21,1]=60029; $data[21,2]=111; $data[21,3]=16589; $data[21,4]=33785; $data[21,5]=0; $data[21,6]=5; $data[21,7]=9655
$data[22,0]='Bielorrusia'; $data[22,1]=53973; $data[22,2]=732; $data[22,3]=30103; $data[22,4]=23562; $data[22,5]=0; $data[22,6]=5; $data[22,7]=308
$data[23,0]='Suecia'; $data[23,1]=51614; $data[23,2]=38; $data[23,3]=0; $data[23,4]=0; $data[23,5]=0; $data[23,6]=0; $data[23,7]=4874
$data[24,0]='Paises Bajos'; $data[24,1]=48783; $data[24,2]=143; $data[24,3]=0; $data[24,4]=0; $data[24,5]=0; $data[24,6]=2; $data[24,7]=6059
$data[25,0]='Colombia'; $data[25,1]=48746; $data[25,2]=0; $data[25,3]=19426; $data[25,4]=27728; $data[25,5]=0; $data[25,6]=0; $data[25,7]=1592
$data[26,0]='Ecuador'; $data[26,1]=46751; $data[26,2]=395; $data[26,3]=23064; $data[26,4]=19791; $data[26,5]=0; $data[26,6]=22; $data[26,7]=3896
$data[27,0]='Egipto'; $data[27,1]=44598; $data[27,2]=1618; $data[27,3]=11931; $data[27,4]=31092; $data[27,5]=0; $data[27,6]=91; $data[27,7]=1575
$data[28,0]='Emiratos Arabes Unidos'; $data[28,1]=42294; $data[28,2]=304; $data[28,3]=27462; $data[28,4]=14543; $data[28,5]=0; $data[28,6]=1; $data[28,7]=289
$data[29,0]='Singapur'; $data[29,1]=40604; $data[29,2]=407; $data[29,3]=29589; $data[29,4]=10989; $data[29,5]=0; $data[29,6]=0; $data[29,7]=26
$data[30,0]='Indonesia'; $data[30,1]=38277; $data[30,2]=857; $data[30,3]=14531; $data[30,4]=21612; $data[30,5]=0; $data[30,6]=43; $data[30,7]=2134
$data[31,0]='Portugal'; $data[31,1]=36690; $data[31,2]=227; $data[31,3]=22669; $data[31,4]=12504; $data[31,5]=0; $data[31,6]=5; $data[31,7]=1517
$data[32,0]='Kuwait'; $data[32,1]=35920; $data[32,2]=454; $data[32,3]=26759; $data[32,4]=8865; $data[32,5]=0; $data[32,6]=7; $data[32,7]=296
$data[33,0]='Ucrania'; $data[33,1]=31154; $data[33,2]=648; $data[33,3]=14082; $data[33,4]=16183; $data[33,5]=0; $data[33,6]=9; $data[33,7]=889
$data[34,0]='Suiza'; $data[34,1]=31117; $data[34,2]=23; $data[34,3]=28800; $data[34,4]=379; $data[34,5]=0; $data[34,6]=0; $data[34,7]=1938
$data[35,0]='Argentina'; $data[35,1]=30295; $data[35,2]=0; $data[35,3]=9564; $data[35,4]=19912; $data[35,5]=0; $data[35,6]=4; $data[35,7]=819
$data[36,0]='Polonia'; $data[36,1]=29392; $data[36,2]=375; $data[36,3]=14226; $data[36,4]=13919; $data[36,5]=0; $data[36,6]=10; $data[36,7]=1247
$data[37,0]='Filipinas'; $data[37,1]=25930; $data[37,2]=538; $data[37,3]=5954; $data[37,4]=18888; $data[37,5]=0; $data[37,6]=14; $data[37,7]=1088
$data[38,0]='Irlanda'; $data[38,1]=25303; $data[38,2]=8; $data[38,3]=22698; $data[38,4]=899; $data[38,5]=0; $data[38,6]=1; $data[38,7]=1706
$data[39,0]='Afganistan'; $data[39,1]=24766; $data[39,2]=664; $data[39,3]=4725; $data[39,4]=19570; $data[39,5]=0; $data[39,6]=20; $data[39,7]=471
$data[40,0]='Oman'; $data[40,1]=23481; $data[40,2]=1404; $data[40,3]=8454; $data[40,4]=14923; $data[40,5]=0; $data[40,6]=5; $data[40,7]=104
$data[41,0]='Republica Dominicana'; $data[41,1]=22962; $data[41,2]=390; $data[41,3]=13320; $data[41,4]=9050; $data[41,5]=0; $data[41,6]=15; $data[41,7]=592
$data[42,0]='Rumania'; $data[42,1]=21999; $data[42,2]=320; $data[42,3]=15719; $data[42,4]=4870; $data[42,5]=0; $data[42,6]=16; $data[42,7]=1410
$data[43,0]='Irak'; $data[43,1]=20209; $data[43,2]=1259; $data[43,3]=8121; $data[43,4]=11481; $data[43,5]=0; $data[43,6]=58; $data[43,7]=607
$data[44,0]='Panama'; $data[44,1]=20059; $data[44,2]=0; $data[44,3]=13759; $data[44,4]=5871; $data[44,5]=0; $data[44,6]=0; $data[44,7]=429
$data[45,0]='Israel'; $data[45,1]=19055; $data[45,2]=83; $data[45,3]=15375; $data[45,4]=3380; $data[45,5]=0; $data[45,6]=0; $data[45,7]=300
$data[46,0]='Barein'; $data[46,1]=18227; $data[46,2]=514; $data[46,3]=12818; $data[46,4]=5367; $data[46,5]=0; $data[46,6]=5; $data[46,7]=42
$data[47,0]='Bolivia'; $data[47,1]=17842; $data[47,2]=913; $data[47,3]=2768; $data[47,4]=14489; $data[47,5]=0; $data[47,6]=26; $data[47,7]=585
$data[48,0]='Japon'; $data[48,1]=17382; $data[48,2]=0; $data[48,3]=15580; $data[48,4]=878; $data[48,5]=0; $data[48,6]=0; $data[48,7]=924
$data[49,0]='Austria'; $data[49,1]=17109; $data[49,2]=31; $data[49,3]=16059; $data[49,4]=373; $data[49,5]=0; $data[49,6]=0; $data[49,7]=677
$data[50,0]='Armenia'; $data[50,1]=16667; $data[50,2]=663; $data[50,3]=6214; $data[50,4]=10184; $data[50,5]=0; $data[50,6]=5; $data[50,7]=269
$data[51,0]='Nigeria'; $data[51,1]=15682; $data[51,2]=0; $data[51,3]=5101; $data[51,4]=10174; $data[51,5]=0; $data[51,6]=0; $data[51,7]=407
$data[52,0]='Kazajistan'; $data[52,1]=14496; $data[52,2]=258; $data[52,3]=9188; $data[52,4]=5231; $data[52,5]=0; $data[52,6]=4; $data[52,7]=77
$data[53,0]='Serbia'; $data[53,1]=12310; $data[53,2]=59; $data[53,3]=11511; $data[53,4]=545; $data[53,5]=0; $data[53,6]=1; $data[53,7]=254
$data[54,0]='Dinamarca'; $data[54,1]=12193; $data[54,2]=54; $data[54,3]=11068; $data[54,4]=528; $data[54,5]=0; $data[54,6]=0; $data[54,7]=597
$data[55,0]='Corea del Sur'; $data[55,1]=12085; $data[55,2]=34; $data[55,3]=10718; $data[55,4]=1090; $data[55,5]=0; $data[55,6]=0; $data[55,7]=277
$data[56,0]='Moldavia'; $data[56,1]=11740; $data[56,2]=281; $data[56,3]=6623; $data[56,4]=4711; $data[56,5]=0; $data[56,6]=8; $data[56,7]=406
$data[57,0]='Ghana'; $data[57,1]=11422; $data[57,2]=304; $data[57,3]=4156; $data[57,4]=7215; $data[57,5]=0; $data[57,6]=3; $data[57,7]=51
$data[58,0]='Argelia'; $data[58,1]=10919; $data[58,2]=109; $data[58,3]=7606; $data[58,4]=2546; $data[58,5]=0; $data[58,6]=7; $data[58,7]=767
$data[59,0]='Chequia'; $data[59,1]=9999; $data[59,2]=8; $data[59,3]=7219; $data[59,4]=2451; $data[59,5]=0; $data[59,6]=1; $data[59,7]=329
$data[60,0]='Azerbaiyan'; $data[60,1]=9957; $data[60,2]=387; $data[60,3]=5583; $data[60,4]=4255; $data[60,5]=0; $data[60,6]=4; $data[60,7]=119
$data[61,0]='Guatemala'; $data[61,1]=9491; $data[61,2]=509; $data[61,3]=1804; $data[61,4]=7320; $data[61,5]=0; $data[61,6]=16; $data[61,7]=367
$data[62,0]='Marruecos'; $data[62,1]=8793; $data[62,2]=101; $data[62,3]=7765; $data[62,4]=816; $data[62,5]=0; $data[62,6]=0; $data[62,7]=212
$data[63,0]='Camerun'; $data[63,1]=8681; $data[63,2]=0; $data[63,3]=4836; $data[63,4]=3633; $data[63,5]=0; $data[63,6]=0; $data[63,7]=212
$data[64,0]='Noruega'; $data[64,1]=8629; $data[64,2]=1; $data[64,3]=8138; $data[64,4]=249; $data[64,5]=0; $data[64,6]=0; $data[64,7]=242
$data[65,0]='Honduras'; $data[65,1]=8455; $data[65,2]=323; $data[65,3]=894; $data[65,4]=7251; $data[65,5]=0; $data[65,6]=4; $data[65,7]=310
$data[66,0]='Malasia'; $data[66,1]=8453; $data[66,2]=8; $data[66,3]=7346; $data[66,4]=986; $data[66,5]=0; $data[66,6]=1; $data[66,7]=121
$data[67,0]='Australia'; $data[67,1]=7320; $data[67,2]=18; $data[67,3]=6838; $data[67,4]=380; $data[67,5]=0; $data[67,6]=0; $data[67,7]=102
$data[68,0]='Finlandia'; $data[68,1]=7104; $data[68,2]=17; $data[68,3]=6200; $data[68,4]=578; $data[68,5]=0; $data[68,6]=1; $data[68,7]=326
$data[69,0]='Sudan'; $data[69,1]=7007; $data[69,2]=0; $data[69,3]=2556; $data[69,4]=4004; $data[69,5]=0; $data[69,6]=0; $data[69,7]=447
$data[70,0]='Nepal'; $data[70,1]=5760; $data[70,2]=425; $data[70,3]=974; $data[70,4]=4767; $data[70,5]=0; $data[70,6]=1; $data[70,7]=19
$data[71,0]='Senegal'; $data[71,1]=5090; $data[71,2]=94; $data[71,3]=3344; $data[71,4]=1686; $data[71,5]=0; $data[71,6]=0; $data[71,7]=60
$data[72,0]='Uzbekistan'; $data[72,1]=5080; $data[72,2]=114; $data[72,3]=3943; $data[72,4]=1118; $data[72,5]=0; $data[72,6]=0; $data[72,7]=19
$data[73,0]='Tayikistan'; $data[73,1]=5035; $data[73,2]=64; $data[73,3]=3409; $data[73,4]=1576; $data[73,5]=0; $data[73,6]=0; $data[73,7]=50
$data[74,0]='Costa de Marfil'; $data[74,1]=4848; $data[74,2]=0; $data[74,3]=2397; $data[74,4]=2406; $data[74,5]=0; $data[74,6]=0; $data[74,7]=45
$data[75,0]='Consejo Danes para los Refugiados'; $data[75,1]=4778; $data[75,2]=54; $data[75,3]=600; $data[75,4]=4071; $data[75,5]=0; $data[75,6]=1; $data[75,7]=107
$data[76,0]='Guinea'; $data[76,1]=4484; $data[76,2]=0; $data[76,3]=3213; $data[76,4]=1246; $data[76,5]=0; $data[76,6]=0; $data[76,7]=25
$data[77,0]='Republica de Yibuti'; $data[77,1]=4465; $data[77,2]=16; $data[77,3]=2950; $data[77,4]=1472; $data[77,5]=0; $data[77,6]=2; $data[77,7]=43
$data[78,0]='Haiti'; $data[78,1]=4165; $data[78,2]=224; $data[78,3]=24; $data[78,4]=4071; $data[78,5]=0; $data[78,6]=6; $data[78,7]=70
$data[79,0]='Luxemburgo'; $data[79,1]=4070; $data[79,2]=7; $data[79,3]=3929; $data[79,4]=31; $data[79,5]=0; $data[79,6]=0; $data[79,7]=110
$data[80,0]='Hungria'; $data[80,1]=4069; $data[80,2]=5; $data[80,3]=2482; $data[80,4]=1025; $data[80,5]=0; $data[80,6]=3; $data[80,7]=562
$data[81,0]='Republica de Macedonia'; $data[81,1]=4057; $data[81,2]=162; $data[81,3]=1710; $data[81,4]=2159; $data[81,5]=0; $data[81,6]=9; $data[81,7]=188
$data[82,0]='El Salvador'; $data[82,1]=3720; $data[82,2]=117; $data[82,3]=1837; $data[82,4]=1811; $data[82,5]=0; $data[82,6]=0; $data[82,7]=72
$data[83,0]='Kenia'; $data[83,1]=3594; $data[83,2]=137; $data[83,3]=1253; $data[83,4]=2238; $data[83,5]=0; $data[83,6]=3; $data[83,7]=103
$data[84,0]='Gabon'; $data[84,1]=3463; $data[84,2]=0; $data[84,3]=1024; $data[84,4]=2416; $data[84,5]=0; $data[84,6]=0; $data[84,7]=23
$data[85,0]='Etiopia'; $data[85,1]=3345; $data[85,2]=179; $data[85,3]=545; $data[85,4]=2743; $data[85,5]=0; $data[85,6]=2; $data[85,7]=57
$data[86,0]='Bulgaria'; $data[86,1]=3266; $data[86,2]=0; $data[86,3]=1723; $data[86,4]=1371; $data[86,5]=0; $data[86,6]=0; $data[86,7]=172
$data[87,0]='Tailandia'; $data[87,1]=3135; $data[87,2]=1; $data[87,3]=2987; $data[87,4]=90; $data[87,5]=0; $data[87,6]=0; $data[87,7]=58
$data[88,0]='Grecia'; $data[88,1]=3121; $data[88,2]=9; $data[88,3]=1374; $data[88,4]=1564; $data[88,5]=0; $data[88,6]=0; $data[88,7]=183
$data[89,0]='Venezuela'; $data[89,1]=2904; $data[89,2]=0; $data[89,3]=487; $data[89,4]=2393; $data[89,5]=0; $data[89,6]=0; $data[89,7]=24
$data[90,0]='Bosnia y Herzegovina'; $data[90,1]=2893; $data[90,2]=0; $data[90,3]=2119; $data[90,4]=611; $data[90,5]=0; $data[90,6]=0; $data[90,7]=163
$data[91,0]='Somalia'; $data[91,1]=2618; $data[91,2]=39; $data[91,3]=577; $data[91,4]=1953; $data[91,5]=0; $data[91,6]=1; $data[91,7]=88
$data[92,0]='Kirguistan'; $data[92,1]=2285; $data[92,2]=78; $data[92,3]=1791; $data[92,4]=467; $data[92,5]=0; $data[92,6]=0; $data[92,7]=27
$data[93,0]='Mayotte'; $data[93,1]=2282; $data[93,2]=0; $data[93,3]=1790; $data[93,4]=464; $data[93,5]=0; $data[93,6]=0; $data[93,7]=28
$data[94,0]='Croacia'; $data[94,1]=2252; $data[94,2]=1; $data[94,3]=2134; $data[94,4]=11; $data[94,5]=0; $data[94,6]=0; $data[94,7]=107
$data[95,0]='Cuba'; $data[95,1]=2248; $data[95,2]=10; $data[95,3]=1948; $data[95,4]=216; $data[95,5]=0; $data[95,6]=0; $data[95,7]=84
$data[96,0]='Republica de Africa Central'; $data[96,1]=2057; $data[96,2]=0; $data[96,3]=363; $data[96,4]=1687; $data[96,5]=0; $data[96,6]=0; $data[96,7]=7
$data[97,0]='Maldivas'; $data[97,1]=2035; $data[97,2]=22; $data[97,3]=1311; $data[97,4]=716; $data[97,5]=0; $data[97,6]=0; $data[97,7]=8
$data[98,0]='Estonia'; $data[98,1]=1973; $data[98,2]=0; $data[98,3]=1705; $data[98,4]=199; $data[98,5]=0; $data[98,6]=0; $data[98,7]=69
$data[99,0]='Sri Lanka'; $data[99,1]=1889; $data[99,2]=5; $data[99,3]=1287; $data[99,4]=591; $data[99,5]=0; $data[99,6]=0; $data[99,7]=11
$data[100,0]='Islandia'; $data[100,1]=1810; $data[100,2]=2; $data[100,3]=1796; $data[100,4]=4; $data[100,5]=0; $data[100,6]=0; $data[100,7]=10
$data[101,0]='Mali'; $data[101,1]=1809; $data[101,2]=33; $data[101,3]=1088; $data[101,4]=617; $data[101,5]=0; $data[101,6]=0; $data[101,7]=104
$data[102,0]='Lituania'; $data[102,1]=1768; $data[102,2]=5; $data[102,3]=1427; $data[102,4]=266; $data[102,5]=0; $data[102,6]=0; $data[102,7]=75
$data[103,0]='Costa Rica'; $data[103,1]=1715; $data[103,2]=53; $data[103,3]=752; $data[103,4]=951; $data[103,5]=0; $data[103,6]=0; $data[103,7]=12
$data[104,0]='Sudan del Sur'; $data[104,1]=1693; $data[104,2]=0; $data[104,3]=49; $data[104,4]=1617; $data[104,5]=0; $data[104,6]=0; $data[104,7]=27
$data[105,0]='Mauritania'; $data[105,1]=1682; $data[105,2]=0; $data[105,3]=311; $data[105,4]=1288; $data[105,5]=0; $data[105,6]=0; $data[105,7]=83
$data[106,0]='Eslovaquia'; $data[106,1]=1548; $data[106,2]=3; $data[106,3]=1410; $data[106,4]=110; $data[106,5]=0; $data[106,6]=0; $data[106,7]=28
$data[107,0]='Albania'; $data[107,1]=1521; $data[107,2]=57; $data[107,3]=1044; $data[107,4]=441; $data[107,5]=0; $data[107,6]=0; $data[107,7]=36
$data[108,0]='Nueva Zelanda'; $data[108,1]=1504; $data[108,2]=0; $data[108,3]=1482; $data[108,4]=0; $data[108,5]=0; $data[108,6]=0; $data[108,7]=22
$data[109,0]='Eslovenia'; $data[109,1]=1495; $data[109,2]=3; $data[109,3]=1359; $data[109,4]=27; $data[109,5]=0; $data[109,6]=0; $data[109,7]=109
$data[110,0]='Nicaragua'; $data[110,1]=1464; $data[110,2]=0; $data[110,3]=953; $data[110,4]=456; $data[110,5]=0; $data[110,6]=0; $data[110,7]=55
$data[111,0]='Guinea-Bisau'; $data[111,1]=1460; $data[111,2]=0; $data[111,3]=153; $data[111,4]=1292; $data[111,5]=0; $data[111,6]=0; $data[111,7]=15
$data[112,0]='Libano'; $data[112,1]=1446; $data[112,2]=4; $data[112,3]=868; $data[112,4]=546; $data[112,5]=0; $data[112,6]=0; $data[112,7]=32
$data[113,0]='Zambia'; $data[113,1]=1358; $data[113,2]=1; $data[113,3]=1122; $data[113,4]=225; $data[113,5]=0; $data[113,6]=1; $data[113,7]=11
$data[114,0]='Guinea Ecuatorial'; $data[114,1]=1306; $data[114,2]=0; $data[114,3]=200; $data[114,4]=1094; $data[114,5]=0; $data[114,6]=0; $data[114,7]=12
$data[115,0]='Paraguay'; $data[115,1]=1289; $data[115,2]=28; $data[115,3]=650; $data[115,4]=628; $data[115,5]=0; $data[115,6]=0; $data[115,7]=11
$data[116,0]='Madagascar'; $data[116,1]=1272; $data[116,2]=20; $data[116,3]=367; $data[116,4]=895; $data[116,5]=0; $data[116,6]=0; $data[116,7]=10
$data[117,0]='Guayana Francesa'; $data[117,1]=1255; $data[117,2]=94; $data[117,3]=534; $data[117,4]=718; $data[117,5]=0; $data[117,6]=1; $data[117,7]=3
$data[118,0]='Sierra Leona'; $data[118,1]=1169; $data[118,2]=37; $data[118,3]=680; $data[118,4]=438; $data[118,5]=0; $data[118,6]=0; $data[118,7]=51
$data[119,0]='Hong Kong'; $data[119,1]=1110; $data[119,2]=0; $data[119,3]=1067; $data[119,4]=39; $data[119,5]=0; $data[119,6]=0; $data[119,7]=4
$data[120,0]='Letonia'; $data[120,1]=1097; $data[120,2]=0; $data[120,3]=845; $data[120,4]=224; $data[120,5]=0; $data[120,6]=0; $data[120,7]=28
$data[121,0]='Tunez'; $data[121,1]=1096; $data[121,2]=2; $data[121,3]=998; $data[121,4]=49; $data[121,5]=0; $data[121,6]=0; $data[121,7]=49
$data[122,0]='Republica de Chipre'; $data[122,1]=980; $data[122,2]=0; $data[122,3]=807; $data[122,4]=155; $data[122,5]=0; $data[122,6]=0; $data[122,7]=18
$data[123,0]='Niger'; $data[123,1]=980; $data[123,2]=0; $data[123,3]=881; $data[123,4]=33; $data[123,5]=0; $data[123,6]=0; $data[123,7]=66
$data[124,0]='Jordania'; $data[124,1]=961; $data[124,2]=8; $data[124,3]=682; $data[124,4]=270; $data[124,5]=0; $data[124,6]=0; $data[124,7]=9
$data[125,0]='Burkina Faso'; $data[125,1]=894; $data[125,2]=2; $data[125,3]=799; $data[125,4]=42; $data[125,5]=0; $data[125,6]=0; $data[125,7]=53
$data[126,0]='Georgia'; $data[126,1]=864; $data[126,2]=13; $data[126,3]=703; $data[126,4]=147; $data[126,5]=0; $data[126,6]=0; $data[126,7]=14
$data[127,0]='Principado de Andorra'; $data[127,1]=853; $data[127,2]=0; $data[127,3]=781; $data[127,4]=21; $data[127,5]=0; $data[127,6]=0; $data[127,7]=51
$data[128,0]='Republica del Chad'; $data[128,1]=850; $data[128,2]=2; $data[128,3]=720; $data[128,4]=57; $data[128,5]=0; $data[128,6]=1; $data[128,7]=73
$data[129,0]='Uruguay'; $data[129,1]=847; $data[129,2]=0; $data[129,3]=784; $data[129,4]=40; $data[129,5]=0; $data[129,6]=0; $data[129,7]=23
$data[130,0]='Cabo Verde'; $data[130,1]=750; $data[130,2]=24; $data[130,3]=301; $data[130,4]=443; $data[130,5]=0; $data[130,6]=0; $data[130,7]=6
$data[131,0]='Yemen'; $data[131,1]=728; $data[131,2]=23; $data[131,3]=53; $data[131,4]=511; $data[131,5]=0; $data[131,6]=4; $data[131,7]=164
$data[132,0]='Congo'; $data[132,1]=728; $data[132,2]=0; $data[132,3]=221; $data[132,4]=483; $data[132,5]=0; $data[132,6]=0; $data[132,7]=24
$data[133,0]='Crucero'; $data[133,1]=712; $data[133,2]=0; $data[133,3]=651; $data[133,4]=48; $data[133,5]=0; $data[133,6]=0; $data[133,7]=13
$data[134,0]='Uganda'; $data[134,1]=696; $data[134,2]=2; $data[134,3]=240; $data[134,4]=456; $data[134,5]=0; $data[134,6]=0; $data[134,7]=0
$data[135,0]='San Marino'; $data[135,1]=694; $data[135,2]=0; $data[135,3]=520; $data[135,4]=132; $data[135,5]=0; $data[135,6]=0; $data[135,7]=42
$data[136,0]='Santo Tome y Principe'; $data[136,1]=659; $data[136,2]=0; $data[136,3]=176; $data[136,4]=471; $data[136,5]=0; $data[136,6]=0; $data[136,7]=12
$data[137,0]='Malta'; $data[137,1]=649; $data[137,2]=3; $data[137,3]=603; $data[137,4]=37; $data[137,5]=0; $data[137,6]=0; $data[137,7]=9
$data[138,0]='Jamaica'; $data[138,1]=615; $data[138,2]=1; $data[138,3]=420; $data[138,4]=185; $data[138,5]=0; $data[138,6]=0; $data[138,7]=10
$data[139,0]='Mozambique'; $data[139,1]=583; $data[139,2]=30; $data[139,3]=151; $data[139,4]=429; $data[139,5]=0; $data[139,6]=1; $data[139,7]=3
$data[140,0]='Ruanda'; $data[140,1]=582; $data[140,2]=41; $data[140,3]=332; $data[140,4]=248; $data[140,5]=0; $data[140,6]=0; $data[140,7]=2
$data[141,0]='Malaui'; $data[141,1]=547; $data[141,2]=18; $data[141,3]=69; $data[141,4]=472; $data[141,5]=0; $data[141,6]=1; $data[141,7]=6
$data[142,0]='Togo'; $data[142,1]=530; $data[142,2]=0; $data[142,3]=291; $data[142,4]=226; $data[142,5]=0; $data[142,6]=0; $data[142,7]=13
$data[143,0]='Tanzania'; $data[143,1]=509; $data[143,2]=0; $data[143,3]=183; $data[143,4]=305; $data[143,5]=0; $data[143,6]=0; $data[143,7]=21
$data[144,0]='Reunion'; $data[144,1]=495; $data[144,2]=6; $data[144,3]=460; $data[144,4]=34; $data[144,5]=0; $data[144,6]=0; $data[144,7]=1
$data[145,0]='Estado de Palestina'; $data[145,1]=492; $data[145,2]=3; $data[145,3]=415; $data[145,4]=74; $data[145,5]=0; $data[145,6]=0; $data[145,7]=3
$data[146,0]='Suazilandia'; $data[146,1]=490; $data[146,2]=4; $data[146,3]=249; $data[146,4]=237; $data[146,5]=0; $data[146,6]=1; $data[146,7]=4
$data[147,0]='Liberia'; $data[147,1]=458; $data[147,2]=12; $data[147,3]=219; $data[147,4]=207; $data[147,5]=0; $data[147,6]=0; $data[147,7]=32
$data[148,0]='Taiwan'; $data[148,1]=443; $data[148,2]=0; $data[148,3]=431; $data[148,4]=5; $data[148,5]=0; $data[148,6]=0; $data[148,7]=7
$data[149,0]='Benin'; $data[149,1]=442; $data[149,2]=30; $data[149,3]=228; $data[149,4]=208; $data[149,5]=0; $data[149,6]=0; $data[149,7]=6
$data[150,0]='Libia'; $data[150,1]=418; $data[150,2]=0; $data[150,3]=62; $data[150,4]=348; $data[150,5]=0; $data[150,6]=0; $data[150,7]=8
$data[151,0]='Zimbabue'; $data[151,1]=383; $data[151,2]=27; $data[151,3]=54; $data[151,4]=325; $data[151,5]=0; $data[151,6]=0; $data[151,7]=4
$data[152,0]='Mauricio'; $data[152,1]=337; $data[152,2]=0; $data[152,3]=325; $data[152,4]=2; $data[152,5]=0; $data[152,6]=0; $data[152,7]=10
$data[153,0]='Isla de Man'; $data[153,1]=336; $data[153,2]=0; $data[153,3]=312; $data[153,4]=0; $data[153,5]=0; $data[153,6]=0; $data[153,7]=24
$data[154,0]='Vietnam'; $data[154,1]=334; $data[154,2]=0; $data[154,3]=323; $data[154,4]=11; $data[154,5]=0; $data[154,6]=0; $data[154,7]=0
$data[155,0]='Montenegro'; $data[155,1]=324; $data[155,2]=0; $data[155,3]=315; $data[155,4]=0; $data[155,5]=0; $data[155,6]=0; $data[155,7]=9
$data[156,0]='Birmania'; $data[156,1]=261; $data[156,2]=0; $data[156,3]=167; $data[156,4]=88; $data[156,5]=0; $data[156,6]=0; $data[156,7]=6
$data[157,0]='Martinica'; $data[157,1]=202; $data[157,2]=0; $data[157,3]=98; $data[157,4]=90; $data[157,5]=0; $data[157,6]=0; $data[157,7]=14
$data[158,0]='Mongolia'; $data[158,1]=197; $data[158,2]=0; $data[158,3]=98; $data[158,4]=99; $data[158,5]=0; $data[158,6]=0; $data[158,7]=0
$data[159,0]='Surinam'; $data[159,1]=196; $data[159,2]=0; $data[159,3]=9; $data[159,4]=184; $data[159,5]=0; $data[159,6]=0; $data[159,7]=3
$data[160,0]='Islas Caimanes'; $data[160,1]=187; $data[160,2]=0; $data[160,3]=115; $data[160,4]=71; $data[160,5]=0; $data[160,6]=0; $data[160,7]=1
$data[161,0]='Islas Feroe'; $data[161,1]=187; $data[161,2]=0; $data[161,3]=187; $data[161,4]=0; $data[161,5]=0; $data[161,6]=0; $data[161,7]=0
$data[162,0]='Siria'; $data[162,1]=177; $data[162,2]=7; $data[162,3]=74; $data[162,4]=97; $data[162,5]=0; $data[162,6]=0; $data[162,7]=6
$data[163,0]='Comoras'; $data[163,1]=176; $data[163,2]=0; $data[163,3]=114; $data[163,4]=60; $data[163,5]=0; $data[163,6]=0; $data[163,7]=2
$data[164,0]='Gibraltar'; $data[164,1]=176; $data[164,2]=0; $data[164,3]=173; $data[164,4]=3; $data[164,5]=0; $data[164,6]=0; $data[164,7]=0
$data[165,0]='Guadalupe'; $data[165,1]=171; $data[165,2]=0; $data[165,3]=157; $data[165,4]=0; $data[165,5]=0; $data[165,6]=0; $data[165,7]=14
$data[166,0]='Guyana'; $data[166,1]=159; $data[166,2]=0; $data[166,3]=95; $data[166,4]=52; $data[166,5]=0; $data[166,6]=0; $data[166,7]=12
$data[167,0]='Bermudas'; $data[167,1]=142; $data[167,2]=0; $data[167,3]=127; $data[167,4]=6; $data[167,5]=0; $data[167,6]=0; $data[167,7]=9
$data[168,0]='Brunei'; $data[168,1]=141; $data[168,2]=0; $data[168,3]=138; $data[168,4]=1; $data[168,5]=0; $data[168,6]=0; $data[168,7]=2
$data[169,0]='Angola'; $data[169,1]=138; $data[169,2]=0; $data[169,3]=61; $data[169,4]=71; $data[169,5]=0; $data[169,6]=0; $data[169,7]=6
$data[170,0]='Camboya'; $data[170,1]=128; $data[170,2]=2; $data[170,3]=125; $data[170,4]=3; $data[170,5]=0; $data[170,6]=0; $data[170,7]=0
$data[171,0]='Trinidad yTobago'; $data[171,1]=118; $data[171,2]=1; $data[171,3]=109; $data[171,4]=1; $data[171,5]=0; $data[171,6]=0; $data[171,7]=8
$data[172,0]='Bahamas'; $data[172,1]=103; $data[172,2]=0; $data[172,3]=68; $data[172,4]=24; $data[172,5]=0; $data[172,6]=0; $data[172,7]=11
$data[173,0]='Aruba'; $data[173,1]=101; $data[173,2]=0; $data[173,3]=98; $data[173,4]=0; $data[173,5]=0; $data[173,6]=0; $data[173,7]=3
$data[174,0]='Monaco'; $data[174,1]=99; $data[174,2]=0; $data[174,3]=93; $data[174,4]=2; $data[174,5]=0; $data[174,6]=0; $data[174,7]=4
$data[175,0]='Eritrea'; $data[175,1]=96; $data[175,2]=31; $data[175,3]=39; $data[175,4]=57; $data[175,5]=0; $data[175,6]=0; $data[175,7]=0
$data[176,0]='Barbados'; $data[176,1]=96; $data[176,2]=0; $data[176,3]=83; $data[176,4]=6; $data[176,5]=0; $data[176,6]=0; $data[176,7]=7
$data[177,0]='Burundi'; $data[177,1]=85; $data[177,2]=0; $data[177,3]=45; $data[177,4]=39; $data[177,5]=0; $data[177,6]=0; $data[177,7]=1
$data[178,0]='Liechtenstein'; $data[178,1]=82; $data[178,2]=0; $data[178,3]=55; $data[178,4]=26; $data[178,5]=0; $data[178,6]=0; $data[178,7]=1
$data[179,0]='San Martin (Parte Holandesa)'; $data[179,1]=77; $data[179,2]=0; $data[179,3]=61; $data[179,4]=1; $data[179,5]=0; $data[179,6]=0; $data[179,7]=15
$data[180,0]='Butan'; $data[180,1]=66; $data[180,2]=4; $data[180,3]=21; $data[180,4]=45; $data[180,5]=0; $data[180,6]=0; $data[180,7]=0
$data[181,0]='Botsuana'; $data[181,1]=60; $data[181,2]=0; $data[181,3]=24; $data[181,4]=35; $data[181,5]=0; $data[181,6]=0; $data[181,7]=1
$data[182,0]='Polinesia Francesa'; $data[182,1]=60; $data[182,2]=0; $data[182,3]=60; $data[182,4]=0; $data[182,5]=0; $data[182,6]=0; $data[182,7]=0
$data[183,0]='Macao'; $data[183,1]=45; $data[183,2]=0; $data[183,3]=45; $data[183,4]=0; $data[183,5]=0; $data[183,6]=0; $data[183,7]=0
$data[184,0]='San Martin (Parte Francesa)'; $data[184,1]=42; $data[184,2]=0; $data[184,3]=36; $data[184,4]=3; $data[184,5]=0; $data[184,6]=0; $data[184,7]=3
$data[185,0]='Puerto Rico'; $data[185,1]=39; $data[185,2]=0; $data[185,3]=1; $data[185,4]=36; $data[185,5]=0; $data[185,6]=0; $data[185,7]=2
$data[186,0]='Guam'; $data[186,1]=32; $data[186,2]=0; $data[186,3]=0; $data[186,4]=31; $data[186,5]=0; $data[186,6]=0; $data[186,7]=1
$data[187,0]='Namibia'; $data[187,1]=32; $data[187,2]=0; $data[187,3]=17; $data[187,4]=15; $data[187,5]=0; $data[187,6]=0; $data[187,7]=0
$data[188,0]='Gambia'; $data[188,1]=28; $data[188,2]=0; $data[188,3]=24; $data[188,4]=3; $data[188,5]=0; $data[188,6]=0; $data[188,7]=1
$data[189,0]='San Vicente y las Granadinas'; $data[189,1]=27; $data[189,2]=0; $data[189,3]=25; $data[189,4]=2; $data[189,5]=0; $data[189,6]=0; $data[189,7]=0
$data[190,0]='Antigua y Barbuda'; $data[190,1]=26; $data[190,2]=0; $data[190,3]=20; $data[190,4]=3; $data[190,5]=0; $data[190,6]=0; $data[190,7]=3
$data[191,0]='Timor Oriental'; $data[191,1]=24; $data[191,2]=0; $data[191,3]=24; $data[191,4]=0; $data[191,5]=0; $data[191,6]=0; $data[191,7]=0
$data[192,0]='Granada'; $data[192,1]=23; $data[192,2]=0; $data[192,3]=22; $data[192,4]=1; $data[192,5]=0; $data[192,6]=0; $data[192,7]=0
$data[193,0]='Curazao'; $data[193,1]=22; $data[193,2]=0; $data[193,3]=15; $data[193,4]=6; $data[193,5]=0; $data[193,6]=0; $data[193,7]=1
$data[194,0]='Nueva Caledonia'; $data[194,1]=21; $data[194,2]=0; $data[194,3]=20; $data[194,4]=1; $data[194,5]=0; $data[194,6]=0; $data[194,7]=0
$data[195,0]='Belice'; $data[195,1]=20; $data[195,2]=0; $data[195,3]=16; $data[195,4]=2; $data[195,5]=0; $data[195,6]=0; $data[195,7]=2
$data[196,0]='Santa Lucia'; $data[196,1]=19; $data[196,2]=0; $data[196,3]=18; $data[196,4]=1; $data[196,5]=0; $data[196,6]=0; $data[196,7]=0
$data[197,0]='Laos'; $data[197,1]=19; $data[197,2]=0; $data[197,3]=19; $data[197,4]=0; $data[197,5]=0; $data[197,6]=0; $data[197,7]=0
$data[198,0]='Dominica'; $data[198,1]=18; $data[198,2]=0; $data[198,3]=16; $data[198,4]=2; $data[198,5]=0; $data[198,6]=0; $data[198,7]=0
$data[199,0]='Fiyi'; $data[199,1]=18; $data[199,2]=0; $data[199,3]=18; $data[199,4]=0; $data[199,5]=0; $data[199,6]=0; $data[199,7]=0
$data[200,0]='Islas Virgenes de los Estados Unidos'; $data[200,1]=17; $data[200,2]=0; $data[200,3]=0; $data[200,4]=17; $data[200,5]=0; $data[200,6]=0; $data[200,7]=0
$data[201,0]='San Cristobal y Nieves'; $data[201,1]=15; $data[201,2]=0; $data[201,3]=15; $data[201,4]=0; $data[201,5]=0; $data[201,6]=0; $data[201,7]=0
$data[202,0]='Islas Malvinas'; $data[202,1]=13; $data[202,2]=0; $data[202,3]=13; $data[202,4]=0; $data[202,5]=0; $data[202,6]=0; $data[202,7]=0
$data[203,0]='Groenlandia'; $data[203,1]=13; $data[203,2]=0; $data[203,3]=13; $data[203,4]=0; $data[203,5]=0; $data[203,6]=0; $data[203,7]=0
$data[204,0]='Santa Sede'; $data[204,1]=12; $data[204,2]=0; $data[204,3]=12; $data[204,4]=0; $data[204,5]=0; $data[204,6]=0; $data[204,7]=0
$data[205,0]='Islas Turcas y Caicos'; $data[205,1]=12; $data[205,2]=0; $data[205,3]=11; $data[205,4]=0; $data[205,5]=0; $data[205,6]=0; $data[205,7]=1
$data[206,0]='Montserrat'; $data[206,1]=11; $data[206,2]=0; $data[206,3]=10; $data[206,4]=0; $data[206,5]=0; $data[206,6]=0; $data[206,7]=1
$data[207,0]='Seychelles'; $data[207,1]=11; $data[207,2]=0; $data[207,3]=11; $data[207,4]=0; $data[207,5]=0; $data[207,6]=0; $data[207,7]=0
$data[208,0]='Sahara Occidental'; $data[208,1]=9; $data[208,2]=0; $data[208,3]=8; $data[208,4]=0; $data[208,5]=0; $data[208,6]=0; $data[208,7]=1
$data[209,0]='Papua Nueva Guinea'; $data[209,1]=8; $data[209,2]=0; $data[209,3]=8; $data[209,4]=0; $data[209,5]=0; $data[209,6]=0; $data[209,7]=0
$data[210,0]='Islas Virgenes Britanicas'; $data[210,1]=8; $data[210,2]=0; $data[210,3]=7; $data[210,4]=0; $data[210,5]=0; $data[210,6]=0; $data[210,7]=1
$data[211,0]='Bonaire, San Eustaquio y Saba'; $data[211,1]=7; $data[211,2]=0; $data[211,3]=7; $data[211,4]=0; $data[211,5]=0; $data[211,6]=0; $data[211,7]=0
$data[212,0]='San Bartolome'; $data[212,1]=6; $data[212,2]=0; $data[212,3]=6; $data[212,4]=0; $data[212,5]=0; $data[212,6]=0; $data[212,7]=0
$data[213,0]='Lesoto'; $data[213,1]=4; $data[213,2]=0; $data[213,3]=2; $data[213,4]=2; $data[213,5]=0; $data[213,6]=0; $data[213,7]=0
$data[214,0]='Anguila'; $data[214,1]=3; $data[214,2]=0; $data[214,3]=3; $data[214,4]=0; $data[214,5]=0; $data[214,6]=0; $data[214,7]=0
$data[215,0]='San Pedro y Miquelon'; $data[215,1]=1; $data[215,2]=0; $data[215,3]=1; $data[215,4]=0; $data[215,5]=0; $data[215,6]=0; $data[215,7]=0

$ws.Range("A4:H219").Value = $data

$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 22:48"

Write-Output "Updated $($data.GetLength(0)) country rows and refreshed the timestamp."
